# Update cryptos list with latest scraped values (GitHub Actions run)
# Note: Price (column D) values are entered with a leading apostrophe so
# Excel stores them as text (matching the source data, which mixes
# thousand-separator formatted numbers like "43.642.44" with plain
# decimals like "7.80" that must keep their trailing zero).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 48 / Row 49 swapped places (ranking changed) plus new values ---
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").Value = "'1.19"
$ws.Range("E48").Value = "  +0.20%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'96.55"
$ws.Range("E49").Value = "  -3.89%  "

# --- Remaining Price (D) / Volume(1h) (E) updates ---
$ws.Range("D2").Value  = "'43.642.44"
$ws.Range("E2").Value  = "  -0.40%  "

$ws.Range("D3").Value  = "'2.283.96"
$ws.Range("E3").Value  = "  -0.73%  "

$ws.Range("E4").Value  = "  +0.07%  "

$ws.Range("D5").Value  = "'96.17"
$ws.Range("E5").Value  = "  +1.91%  "

$ws.Range("D6").Value  = "'267.05"
$ws.Range("E6").Value  = "  -0.75%  "

$ws.Range("D7").Value  = "'0.622"
$ws.Range("E7").Value  = "  -0.48%  "

$ws.Range("E8").Value  = "  -0.13%  "

$ws.Range("D9").Value  = "'0.609"
$ws.Range("E9").Value  = "  -2.35%  "

$ws.Range("D10").Value = "'45.91"
$ws.Range("E10").Value = "  +0.27%  "

$ws.Range("D11").Value = "'0.0935"
$ws.Range("E11").Value = "  -0.03%  "

$ws.Range("D12").Value = "'7.80"
$ws.Range("E12").Value = "  -2.80%  "

$ws.Range("E13").Value = "  +0.23%  "

$ws.Range("D14").Value = "'2.628.41"
$ws.Range("E14").Value = "  -0.62%  "

$ws.Range("E15").Value = "  -1.86%  "

$ws.Range("D16").Value = "'0.845"
$ws.Range("E16").Value = "  -1.00%  "

$ws.Range("D17").Value = "'2.288.21"
$ws.Range("E17").Value = "  -0.49%  "

$ws.Range("D18").Value = "'43.607.38"
$ws.Range("E18").Value = "  -0.45%  "

$ws.Range("D19").Value = "'0.0000107"
$ws.Range("E19").Value = "  +1.82%  "

$ws.Range("D20").Value = "'6.20"
$ws.Range("E20").Value = "  -1.59%  "

$ws.Range("D21").Value = "'72.13"
$ws.Range("E21").Value = "  +1.33%  "

$ws.Range("D22").Value = "'2.47"
$ws.Range("E22").Value = "  +7.87%  "

$ws.Range("D23").Value = "'232.69"
$ws.Range("E23").Value = "  -1.81%  "

$ws.Range("E24").Value = "  -5.70%  "

$ws.Range("E26").Value = "  +0.75%  "

$ws.Range("E27").Value = "  -1.48%  "

$ws.Range("E28").Value = "  +2.60%  "

$ws.Range("D29").Value = "'40.61"
$ws.Range("E29").Value = "  +3.42%  "

$ws.Range("D30").Value = "'2.27"
$ws.Range("E30").Value = "  +0.89%  "

$ws.Range("D31").Value = "'176.39"
$ws.Range("E31").Value = "  +1.62%  "

$ws.Range("D32").Value = "'21.80"
$ws.Range("E32").Value = "  -2.32%  "

$ws.Range("E33").Value = "  -0.94%  "

$ws.Range("E34").Value = "  -3.62%  "

$ws.Range("E35").Value = "  +0.57%  "

$ws.Range("E36").Value = "  -1.87%  "

$ws.Range("E37").Value = "  +0.71%  "

$ws.Range("D38").Value = "'4.33"
$ws.Range("E38").Value = "  -3.56%  "

$ws.Range("D39").Value = "'3.42"
$ws.Range("E39").Value = "  +0.05%  "

$ws.Range("E40").Value = "  -0.68%  "

$ws.Range("E41").Value = "  -0.34%  "

$ws.Range("D42").Value = "'12.28"
$ws.Range("E42").Value = "  -0.21%  "

$ws.Range("D43").Value = "'1.34"
$ws.Range("E43").Value = "  +0.62%  "

$ws.Range("D44").Value = "'64.88"
$ws.Range("E44").Value = "  +5.55%  "

$ws.Range("E45").Value = "  -0.38%  "

$ws.Range("E46").Value = "  -5.03%  "

$ws.Range("E47").Value = "  -0.27%  "

$ws.Range("D50").Value = "'0.436"
$ws.Range("E50").Value = "  +1.19%  "

$ws.Range("D51").Value = "'2.508.31"
$ws.Range("E51").Value = "  -0.55%  "
